# Applies the "ALL exchange logic added to income statement" update to the
# Issue Tracker workbook: adds new EXCHANGE FILTER ALL / BULK WATERMARK
# columns, new UPDATE RUN TIME / UPDATE DAY / UPDATE TIME / UPDATE CREDIT
# CONSUMPTION columns, fills in data for several workflow rows, and removes
# the old "PULLS ALL QUARTERS" note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("A1").Value = "WORKFLOW"
$ws.Range("B1").Value = "WORKS"
$ws.Range("C1").Value = "PROBLEM"
$ws.Range("D1").Value = "PRIVATE KEY WORKS"
$ws.Range("E1").Value = "EXCHANGE FILTER ALL"
$ws.Range("F1").Value = "BULK WATERMARK"
$ws.Range("G1").Value = "SYMBOL PULL OPTIMIZED?"
$ws.Range("H1").Value = "UPDATE RUN TIME"
$ws.Range("I1").Value = "UPDATE DAY"
$ws.Range("J1").Value = "UPDATE TIME"
$ws.Range("K1").Value = "UPDATE CREDIT CONSUMPTION"

# ---- Row 2: earnings call transcript ----
$ws.Range("A2").Value = "earnings call transcript"
$ws.Range("B2").Value = "YES"
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = "YES"
$ws.Range("E2").Value = "YES"
$ws.Range("F2").Value = "YES"
$ws.Range("G2").Value = "YES"
$ws.Range("H2").Value = 40
$ws.Range("I2").Value = "FRIDAY"
$ws.Range("J2").Value = 0.75
$ws.Range("J2").NumberFormat = "h:mm AM/PM"

# ---- Row 3: balance sheet ----
$ws.Range("A3").Value = "balance sheet"
$ws.Range("B3").Value = "YES"
$ws.Range("D3").Value = "YES"
$ws.Range("E3").Value = "YES"
$ws.Range("F3").Value = "YES"
$ws.Range("G3").Value = "YES"

# ---- Row 4: cash flow ----
$ws.Range("A4").Value = "cash flow"
$ws.Range("B4").Value = "YES"
$ws.Range("D4").Value = "YES"
$ws.Range("E4").Value = "YES"
$ws.Range("F4").Value = "YES"
$ws.Range("G4").Value = "YES"
$ws.Range("H4").Value = 65
$ws.Range("K4").Value = 0.1

# ---- Row 5: income statement ----
$ws.Range("A5").Value = "income statement"

# ---- Row 6: insider transactions ----
$ws.Range("A6").Value = "insider transactions"
$ws.Range("B6").Value = "YES"
$ws.Range("C6").Value = "not printing iters to log"
$ws.Range("D6").Value = "NO"
$ws.Range("F6").Value = "NO"

# ---- Row 7: etf profile ----
$ws.Range("A7").Value = "etf profile"
$ws.Range("B7").Value = "YES"
$ws.Range("I7").Value = "NO SCHEDULE"

# ---- Row 8: time series ----
$ws.Range("A8").Value = "time series"

# ---- Row 9: fred commodities ----
$ws.Range("A9").Value = "fred commodities"
$ws.Range("B9").Value = "YES"
$ws.Range("D9").Value = "YES"
$ws.Range("E9").Value = "N/A"
$ws.Range("F9").Value = "N/A"
$ws.Range("G9").Value = "N/A"
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = "DAILY"
$ws.Range("J9").Value = 15/1440
$ws.Range("J9").NumberFormat = "h:mm AM/PM"

# ---- Row 10: fred econ indicators ----
$ws.Range("A10").Value = "fred econ indicators"
$ws.Range("B10").Value = "YES"
$ws.Range("D10").Value = "YES"
$ws.Range("E10").Value = "N/A"
$ws.Range("F10").Value = "N/A"
$ws.Range("G10").Value = "N/A"
$ws.Range("H10").Value = 10
$ws.Range("I10").Value = "DAILY"
$ws.Range("J10").Value = 0
$ws.Range("J10").NumberFormat = "h:mm AM/PM"

# ---- Column widths (best effort match of original autofit values; the
#      interop layer quantizes width to 1/6-character steps so an exact
#      match to Excel's native 1/256 step width isn't achievable) ----
$ws.Columns.Item(5).ColumnWidth = 18.28
$ws.Columns.Item(6).ColumnWidth = 16.59
$ws.Columns.Item(7).ColumnWidth = 21.78
$ws.Columns.Item(8).ColumnWidth = 14.78

# ---- Selection, matching the saved cursor position in the diff ----
$ws.Range("G21").Select() | Out-Null

Write-Host "Applied ALL exchange logic changes to income statement sheet"
